$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: num_customers (C36) increases from 122 to 123; retention_rate (E36) recalculated
$ws.Range("C36").Value = 123
$ws.Range("E36").Value = 123/1930

# Row 37: num_customers (C37) and cohort_size (D37) increase from 772 to 783; retention_rate (E37) stays 1
$ws.Range("C37").Value = 783
$ws.Range("D37").Value = 783
$ws.Range("E37").Value = 783/783
